{"js": "// Add \"Square One Standard \" before \"Statement of Work\" in the document\n// title (first paragraph), matching the surrounding bold/Times New Roman\n// formatting, per commit message: 'Add \"Square One Standard\" to template\n// titles'.\n\nconst titlePara = context.document.body.paragraphs.getFirst();\nconst results = titlePara.search(\"Statement of Work\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Insert immediately before the matched text; inherits the matched\n  // range's (bold, Times New Roman) character formatting.\n  results.items[0].insertText(\"Square One Standard \", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "# Add \"Square One Standard \" before \"Statement of Work\" in the document\n# title, matching the surrounding bold/Times New Roman formatting, per\n# commit message: 'Add \"Square One Standard\" to template titles'.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Statement of Work\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    # $range now spans the found \"Statement of Work\" text; inserting\n    # immediately before it inherits that text's character formatting.\n    $range.InsertBefore(\"Square One Standard \")\n}\n"}
